$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the specific Price cells we are about to update,
# so values like "1.00" or "254.23" remain text instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '97.490.35'
$ws.Range("E2").Value = '  -1.66%  '

$ws.Range("D3").Value = '3.397.32'
$ws.Range("E3").Value = '  +3.37%  '

$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").Value = '254.23'
$ws.Range("E5").Value = '  -0.05%  '

$ws.Range("D6").Value = '647.65'
$ws.Range("E6").Value = '  +3.88%  '

$ws.Range("E7").Value = '  -1.02%  '

$ws.Range("D8").Value = '0.422'
$ws.Range("E8").Value = '  +5.25%  '

$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("D10").Value = '1.05'
$ws.Range("E10").Value = '  +6.04%  '

$ws.Range("D11").Value = '3.393.12'
$ws.Range("E11").Value = '  +3.35%  '

$ws.Range("E12").Value = '  +4.57%  '

$ws.Range("D13").Value = '41.11'
$ws.Range("E13").Value = '  +3.53%  '

$ws.Range("E14").Value = '  +13.06%  '

$ws.Range("E15").Value = '  +2.54%  '

$ws.Range("D16").Value = '97.169.48'
$ws.Range("E16").Value = '  -1.70%  '

$ws.Range("D17").Value = '4.014.12'
$ws.Range("E17").Value = '  +2.76%  '

$ws.Range("D18").Value = '8.44'
$ws.Range("E18").Value = '  +32.89%  '

$ws.Range("D19").Value = '3.398.97'
$ws.Range("E19").Value = '  +3.39%  '

$ws.Range("D20").Value = '17.26'
$ws.Range("E20").Value = '  +13.00%  '

$ws.Range("D21").Value = '10.68'
$ws.Range("E21").Value = '  +14.64%  '

$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").Value = '3.38'
$ws.Range("E22").Value = '  -2.49%  '

$ws.Range("B23").Value = 'Stellar'
$ws.Range("C23").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D23").Value = '0.478'
$ws.Range("E23").Value = '  +39.79%  '

$ws.Range("D24").Value = '500.01'
$ws.Range("E24").Value = '  +2.09%  '

$ws.Range("D25").Value = '0.0000202'
$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '6.02'
$ws.Range("E26").Value = '  +6.74%  '

$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '95.27'
$ws.Range("E27").Value = '  +6.90%  '

$ws.Range("D28").Value = '12.52'
$ws.Range("E28").Value = '  +3.44%  '

$ws.Range("D29").Value = '3.578.21'
$ws.Range("E29").Value = '  +3.36%  '

$ws.Range("E30").Value = '  +9.14%  '

$ws.Range("D31").Value = '0.196'
$ws.Range("E31").Value = '  +3.56%  '

$ws.Range("E32").Value = '  -0.34%  '

$ws.Range("E33").Value = '  +7.17%  '

$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").Value = '0.563'
$ws.Range("E35").Value = '  +18.36%  '

$ws.Range("D36").Value = '29.35'
$ws.Range("E36").Value = '  +4.68%  '

$ws.Range("E37").Value = '  +14.43%  '

$ws.Range("D38").Value = '7.61'
$ws.Range("E38").Value = '  +5.13%  '

$ws.Range("E39").Value = '  +1.19%  '

$ws.Range("E40").Value = '  +12.36%  '

$ws.Range("D41").Value = '504.90'
$ws.Range("E41").Value = '  +3.46%  '

$ws.Range("D42").Value = '24.69'
$ws.Range("E42").Value = '  -0.25%  '

$ws.Range("D43").Value = '0.851'
$ws.Range("E43").Value = '  +9.95%  '

$ws.Range("D44").Value = '3.64'
$ws.Range("E44").Value = '  -2.67%  '

$ws.Range("D45").Value = '0.0409'
$ws.Range("E45").Value = '  +20.52%  '

$ws.Range("D46").Value = '5.45'
$ws.Range("E46").Value = '  +14.36%  '

$ws.Range("E47").Value = '  -0.04%  '

$ws.Range("D48").Value = '3.18'
$ws.Range("E48").Value = '  +2.48%  '

$ws.Range("D49").Value = '8.11'
$ws.Range("E49").Value = '  +10.83%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '51.43'
$ws.Range("E50").Value = '  +11.19%  '

$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").Value = '1.55'
$ws.Range("E51").Value = '  +13.52%  '
